# Insert a new data row above row 147 (shifting existing rows 147-167 down to 148-168)
# and populate it with a new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(147).Insert()

$ws.Range("A147").Value = 10
$ws.Range("B147").Value = "Vega Modelo de Temuco"
$ws.Range("C147").Value = "La Araucanía"
$ws.Range("D147").Value = 44776
$ws.Range("E147").Value = 9
$ws.Range("F147").Value = 100112012
$ws.Range("G147").Value = "Espinaca"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 50
$ws.Range("K147").Value = 9000
$ws.Range("L147").Value = 9000
$ws.Range("M147").Value = 9000
$ws.Range("N147").Value = "`$/docena de atados"
$ws.Range("O147").Value = "Región de Coquimbo"
$ws.Range("P147").Value = 3000
$ws.Range("Q147").Value = 3
$ws.Range("R147").Value = "Hortaliza"
